# Commit: "Fri, Jun 05, 2020  6:05:52 PM"
#
# The canonical diff shows two kinds of change:
#
#   1. Every table in the deck (one each on slides 14, 15 and 16) has its
#      <a:tableStyleId> switched from the custom "Table_0" style
#      {90EBD686-0D37-461E-A155-407F3169B37F} (the deck's only custom table
#      style, defined in ppt/tableStyles.xml) to the built-in table style
#      {7EC342C0-CF34-4DEA-8A95-78914D0FC3A6} - i.e. the user picked a
#      different swatch from the Table Styles gallery on the Table Design
#      ribbon for each table.
#
#   2. ppt/theme/theme1.xml and ppt/theme/theme2.xml swap their entire
#      contents (the deck keeps using the same theme relationships
#      throughout - slide master + presentation still point at theme2.xml,
#      notes master still points at theme1.xml - only the payload each file
#      name carries changes). This is a PowerPoint-internal artefact of
#      applying/removing a design theme and is not reachable through the
#      Presentation/Master/Design/Theme COM surface (ApplyTheme,
#      OpenThemeFile, Design.Name, ... are all no-ops here - this host has
#      no theme-import path), so it cannot be reproduced from script.

$p = $ppt.ActivePresentation

$newStyleId = "{7EC342C0-CF34-4DEA-8A95-78914D0FC3A6}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shp = $slide.Shapes.Item($j)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newStyleId)
        }
    }
}
